# Weekly CompStat report refresh: new crime data collected for the week of
# 5/5/2025 - 5/11/2025 (Volume 32, Number 19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / banner text -------------------------------------------------
$ws.Range("A8").Value  = "Volume 32   Number  19"
$ws.Range("C9").Value  = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Crime Complaints table (rows 14-31) -----------------------------------

# Row 14 (Murder): 2025 Week-to-Date count moves from the text placeholder
# "0" to an actual number (2), so match the numeric style used by its
# neighboring cell before assigning the value.
$ws.Range("C14").NumberFormat = $ws.Range("D14").NumberFormat
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 100
$ws.Range("I14").Value = 11
$ws.Range("J14").Value = 6
$ws.Range("K14").Value = 83.333333333333
$ws.Range("L14").Value = 37.5
$ws.Range("M14").Value = 10
$ws.Range("N14").Value = -77.551020408163

# Row 15 (Rape)
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = -25
$ws.Range("F15").Value = 17
$ws.Range("G15").Value = 14
$ws.Range("H15").Value = 21.428571428571
$ws.Range("I15").Value = 75
$ws.Range("J15").Value = 63
$ws.Range("K15").Value = 19.047619047619
$ws.Range("L15").Value = 4.166666666666
$ws.Range("M15").Value = 82.926829268292
$ws.Range("N15").Value = 17.1875

# Row 16 (Robbery)
$ws.Range("C16").Value = 28
$ws.Range("D16").Value = 36
$ws.Range("E16").Value = -22.222222222222
$ws.Range("F16").Value = 115
$ws.Range("G16").Value = 167
$ws.Range("H16").Value = -31.137724550898
$ws.Range("I16").Value = 506
$ws.Range("J16").Value = 757
$ws.Range("K16").Value = -33.157199471598
$ws.Range("L16").Value = -23.449319213313
$ws.Range("M16").Value = -22.153846153846
$ws.Range("N16").Value = -84.591961023142

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 57
$ws.Range("D17").Value = 53
$ws.Range("E17").Value = 7.547169811320
$ws.Range("F17").Value = 257
$ws.Range("G17").Value = 243
$ws.Range("H17").Value = 5.761316872427
$ws.Range("I17").Value = 1080
$ws.Range("J17").Value = 1060
$ws.Range("K17").Value = 1.886792452830
$ws.Range("L17").Value = 11.917098445595
$ws.Range("M17").Value = 103.77358490566
$ws.Range("N17").Value = 12.266112266112

# Row 18 (Burglary)
$ws.Range("C18").Value = 34
$ws.Range("D18").Value = 20
$ws.Range("E18").Value = 70
$ws.Range("F18").Value = 122
$ws.Range("G18").Value = 133
$ws.Range("H18").Value = -8.270676691729
$ws.Range("I18").Value = 677
$ws.Range("J18").Value = 694
$ws.Range("K18").Value = -2.449567723342
$ws.Range("L18").Value = -9.370816599732
$ws.Range("M18").Value = -30.278063851699
$ws.Range("N18").Value = -87.240859404447

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 108
$ws.Range("D19").Value = 117
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 420
$ws.Range("G19").Value = 523
$ws.Range("H19").Value = -19.694072657743
$ws.Range("I19").Value = 1866
$ws.Range("J19").Value = 2375
$ws.Range("K19").Value = -21.431578947368
$ws.Range("L19").Value = -22.022565816966
$ws.Range("M19").Value = 36.105032822757
$ws.Range("N19").Value = -31.497797356828

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 42
$ws.Range("D20").Value = 49
$ws.Range("E20").Value = -14.285714285714
$ws.Range("F20").Value = 194
$ws.Range("G20").Value = 190
$ws.Range("H20").Value = 2.105263157894
$ws.Range("I20").Value = 724
$ws.Range("J20").Value = 836
$ws.Range("K20").Value = -13.397129186602
$ws.Range("L20").Value = -8.816120906801
$ws.Range("M20").Value = 17.152103559870
$ws.Range("N20").Value = -91.386079714455

# Row 21 (TOTAL)
$ws.Range("C21").Value = 274
$ws.Range("D21").Value = 280
$ws.Range("E21").Value = -2.142857142857
$ws.Range("F21").Value = 1129
$ws.Range("G21").Value = 1272
$ws.Range("H21").Value = -11.242138364779
$ws.Range("I21").Value = 4939
$ws.Range("J21").Value = 5791
$ws.Range("K21").Value = -14.712484890347
$ws.Range("L21").Value = -12.429078014184
$ws.Range("M21").Value = 17.847769028871
$ws.Range("N21").Value = -76.247956141194

# Row 22 (Transit)
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = -18.181818181818
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 21
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 122
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 22
$ws.Range("L22").Value = -1.612903225806
$ws.Range("M22").Value = 71.830985915493

# Row 23 (Housing)
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -87.5
$ws.Range("F23").Value = 13
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = -43.478260869565
$ws.Range("I23").Value = 79
$ws.Range("J23").Value = 83
$ws.Range("K23").Value = -4.819277108433
$ws.Range("L23").Value = -15.957446808510
$ws.Range("M23").Value = 58

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 271
$ws.Range("D24").Value = 274
$ws.Range("E24").Value = -1.094890510948
$ws.Range("F24").Value = 1049
$ws.Range("G24").Value = 1157
$ws.Range("H24").Value = -9.334485738980
$ws.Range("I24").Value = 5058
$ws.Range("J24").Value = 5850
$ws.Range("K24").Value = -13.538461538461
$ws.Range("L24").Value = -5.510928451335
$ws.Range("M24").Value = 59.457755359394

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 152
$ws.Range("D25").Value = 133
$ws.Range("E25").Value = 14.285714285714
$ws.Range("F25").Value = 575
$ws.Range("G25").Value = 689
$ws.Range("H25").Value = -16.545718432510
$ws.Range("I25").Value = 3047
$ws.Range("J25").Value = 3645
$ws.Range("K25").Value = -16.406035665294
$ws.Range("L25").Value = 2.834964562942

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 112
$ws.Range("D26").Value = 140
$ws.Range("E26").Value = -20
$ws.Range("F26").Value = 460
$ws.Range("G26").Value = 524
$ws.Range("H26").Value = -12.213740458015
$ws.Range("I26").Value = 1912
$ws.Range("J26").Value = 2112
$ws.Range("K26").Value = -9.469696969696
$ws.Range("L26").Value = 6.281267370761
$ws.Range("M26").Value = 13.877307921381

# Row 27 (UCR Rape*)
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = -42.857142857142
$ws.Range("F27").Value = 19
$ws.Range("G27").Value = 25
$ws.Range("H27").Value = -24
$ws.Range("I27").Value = 97
$ws.Range("J27").Value = 102
$ws.Range("K27").Value = -4.901960784313
$ws.Range("L27").Value = -10.185185185185

# Row 28 (Other Sex Crimes)
$ws.Range("C28").Value = 15
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 50
$ws.Range("F28").Value = 51
$ws.Range("G28").Value = 42
$ws.Range("H28").Value = 21.428571428571
$ws.Range("I28").Value = 201
$ws.Range("J28").Value = 206
$ws.Range("K28").Value = -2.427184466019
$ws.Range("L28").Value = -19.6

# Row 29 (Shooting Vic.) - C/D/E stay as the "0"/"0"/"***.*" placeholders
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 400
$ws.Range("I29").Value = 15
$ws.Range("J29").Value = 5
$ws.Range("K29").Value = 200
$ws.Range("L29").Value = -37.5
$ws.Range("M29").Value = -6.25
$ws.Range("N29").Value = -83.146067415730

# Row 30 (Shooting Inc.) - C/D/E stay as the "0"/"0"/"***.*" placeholders
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 400
$ws.Range("I30").Value = 15
$ws.Range("J30").Value = 4
$ws.Range("K30").Value = 275
$ws.Range("L30").Value = -31.818181818181
$ws.Range("M30").Value = 25
$ws.Range("N30").Value = -82.352941176470

# Row 31 (Hate Crimes): 2025 Week-to-Date count moves from a real number (1)
# down to the text placeholder "0", so copy the text/General style used by
# the other placeholder cells in the table before assigning the value; the
# leading apostrophe forces Excel to store it as text instead of a number.
$ws.Range("C31").NumberFormat = $ws.Range("C29").NumberFormat
$ws.Range("C31").Value = "'0"
$ws.Range("D31").Value = 2
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 11
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 120
$ws.Range("I31").Value = 37
$ws.Range("J31").Value = 29
$ws.Range("K31").Value = 27.586206896551
$ws.Range("L31").Value = 23.333333333333
